## properties to properties tabs of tourney sheets
##
## - remove the competition-key / host-key / venue-key.N rows from the
##   Tournament sheet (their data now lives in a dedicated Properties sheet)
## - add a new "Properties" sheet (key/value/notes) at the end of the
##   workbook with that same data, plus color + timezone lookups

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Tournament sheet: drop the rows that held competition-key / host-key
#    (rows 2:3) and the venue-key.1..11 rows (originally rows 18:28, now
#    16:26 once the first two rows are gone).
# ---------------------------------------------------------------------
$tournament = $wb.Worksheets.Item("Tournament")
$tournament.Range("A2:A3").EntireRow.Delete()
$tournament.Range("A16:A26").EntireRow.Delete()
$tournament.Range("A4:B4").Select()

# ---------------------------------------------------------------------
# 2. Colors sheet: nothing data-wise changes here, just where the
#    selection ends up.
# ---------------------------------------------------------------------
$colors = $wb.Worksheets.Item("Colors")
$colors.Columns("K").Select()

# ---------------------------------------------------------------------
# 3. New "Properties" sheet, appended after "#Work".
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$properties = $wb.Worksheets.Add($null, $lastSheet)
$properties.Name = "Properties"

$rows = @(
  @("key", "value", "notes"),
  @("competition", "mens-club-world-cup", ""),
  @("host", "usa", ""),
  @("timezone", "US/Eastern", ""),
  @("color.a", "#c4e1b5", "green"),
  @("color.b", "#fab077", "orange"),
  @("color.c", "#94d9f5", "cyan"),
  @("color.d", "#fee289", "yellow"),
  @("color.e", "#eecbef", "purple"),
  @("color.f", "#f79d8f", "red"),
  @("color.g", "#b0d0ee", "blue"),
  @("color.h", "#eb84af", "rose"),
  @("venue.01", "us-atlanta-ga", ""),
  @("venue.02", "us-charlotte-nc", ""),
  @("venue.03", "us-cincinnati-oh", ""),
  @("venue.04", "us-pasadena-ca", ""),
  @("venue.05", "us-miami-fl", ""),
  @("venue.06", "us-nashville-tn", ""),
  @("venue.07", "us-east-rutherford-nj", ""),
  @("venue.08", "us-orlando-fl", ""),
  @("venue.09", "us-philadelphia-pa", ""),
  @("venue.10", "us-seattle-wa", ""),
  @("venue.11", "us-washington-dc", "")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
  for ($c = 0; $c -lt 3; $c++) {
    $value = $rows[$r][$c]
    if ($value -ne "") {
      $properties.Cells.Item($r + 1, $c + 1).Value = $value
    }
  }
}

$propertiesTable = $properties.ListObjects.Add(1, $properties.Range("A1:C23"), $null, 1)
$propertiesTable.Name = "Properties"
$propertiesTable.TableStyle = "TableStyleMedium2"

$properties.Range("B5").Select()

# ---------------------------------------------------------------------
# 4. Leave the workbook focused back on the Tournament tab, matching the
#    original author's final view.
# ---------------------------------------------------------------------
$tournament.Select()
